# Generate Report for Handback
# This script updates the Overview/zh-cn/de-de "Ready for handoff" status
# to "Handback transform failed" and records the mismatch error detail
# explaining why the handback transform failed for both the zh-cn and
# de-de target files, widening the "Error Detail" column to fit.

$wb = $excel.ActiveWorkbook

$zhMessage = "Handback file name: bweebyik.2ij is different with handoff file name: ad98c24e-8692-4f07-866f-d65eb5bd6eb1.b7000313cce4b1ac24baa2bc4984ec8585d0fade.zh-cn."
$deMessage = "Handback file name: bweebyik.2ij is different with handoff file name: ad98c24e-8692-4f07-866f-d65eb5bd6eb1.b7000313cce4b1ac24baa2bc4984ec8585d0fade.de-de."

# Update status text on every sheet that references it (Overview, zh-cn, de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handback transform failed"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handback transform failed"

# Record the error detail describing the handback/handoff file name mismatch
$wsZh.Range("P3").Value = $zhMessage
$wsDe.Range("P3").Value = $deMessage

# Widen the Error Detail column (column P) on both locale sheets to fit the message.
# Excel stores column widths in "characters of the workbook's Normal font" but rounds
# to whole pixels internally (padding of 5px, MDW=7px for Calibri 11) before
# re-deriving the character width used in the saved XML, so setting ColumnWidth to
# exactly 40 actually persists as ~40.83. Using 39 + 1/7 lands on the same pixel
# width (280px) that yields a persisted width of exactly 40.
$wsZh.Columns.Item(16).ColumnWidth = 39 + 1/7
$wsDe.Columns.Item(16).ColumnWidth = 39 + 1/7
